$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 286 (old rows 286-309 shift down to 288-311)
$ws.Rows.Item(286).Insert()
$ws.Rows.Item(286).Insert()

# New row 286: Acelga, Primera, week of 2022-01-17 (serial 44578)
$ws.Range("A286").Value = 8
$ws.Range("B286").Value = "Terminal La Palmera de La Serena"
$ws.Range("C286").Value = "Coquimbo"
$ws.Range("D286").Value = 44578
$ws.Range("E286").Value = 4
$ws.Range("F286").Value = 100112009
$ws.Range("G286").Value = "Acelga"
$ws.Range("H286").Value = "Sin especificar"
$ws.Range("I286").Value = "Primera"
$ws.Range("J286").Value = 2400
$ws.Range("K286").Value = 450
$ws.Range("L286").Value = 500
$ws.Range("M286").Value = 475
$ws.Range("N286").Value = "$/atado 1,5 a 2 kilos"
$ws.Range("O286").Value = "Provincia del Elquí"
$ws.Range("P286").Value = 238
$ws.Range("Q286").Value = 2
$ws.Range("R286").Value = "Hortaliza"

# New row 287: Acelga, Segunda, week of 2022-01-17 (serial 44578)
$ws.Range("A287").Value = 8
$ws.Range("B287").Value = "Terminal La Palmera de La Serena"
$ws.Range("C287").Value = "Coquimbo"
$ws.Range("D287").Value = 44578
$ws.Range("E287").Value = 4
$ws.Range("F287").Value = 100112009
$ws.Range("G287").Value = "Acelga"
$ws.Range("H287").Value = "Sin especificar"
$ws.Range("I287").Value = "Segunda"
$ws.Range("J287").Value = 1480
$ws.Range("K287").Value = 350
$ws.Range("L287").Value = 400
$ws.Range("M287").Value = 375
$ws.Range("N287").Value = "$/atado 1,5 a 2 kilos"
$ws.Range("O287").Value = "Provincia del Elquí"
$ws.Range("P287").Value = 188
$ws.Range("Q287").Value = 2
$ws.Range("R287").Value = "Hortaliza"
